$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new progress note for Anushka Chincholkar (row 4, column D)
$ws.Range("D4").Value = "Studied basics of Flutter UI Design"

# Match the saved selection state from the diff (active cell D4)
$ws.Range("D4").Select()
